# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# timestamps for the 407bfbc6-... row on the zh-cn and de-de status sheets,
# to reflect a freshly regenerated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-17 04:15:58"
$wsZhCn.Range("G4").Value = "2016-02-17 04:16:43"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-17 04:16:08"
$wsDeDe.Range("G4").Value = "2016-02-17 04:17:02"
